$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "67.734.49"
$ws.Range("E2").Value = "  -4.01%  "

$ws.Range("D3").Value = "3.439.79"
$ws.Range("E3").Value = "  -5.02%  "

Set-TextValue "D5" "580.36"
$ws.Range("E5").Value = "  -0.55%  "

Set-TextValue "D6" "164.64"
$ws.Range("E6").Value = "  -5.98%  "

$ws.Range("E7").Value = "  -4.73%  "

$ws.Range("D8").Value = "3.434.39"
$ws.Range("E8").Value = "  -4.80%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("E10").Value = "  -4.91%  "

$ws.Range("E11").Value = "  +1.52%  "

$ws.Range("E12").Value = "  -7.78%  "

Set-TextValue "D13" "45.84"
$ws.Range("E13").Value = "  -5.12%  "

Set-TextValue "D14" "0.0000268"
$ws.Range("E14").Value = "  -4.47%  "

$ws.Range("D15").Value = "3.997.45"
$ws.Range("E15").Value = "  -5.12%  "

Set-TextValue "D16" "610.86"
$ws.Range("E16").Value = "  -10.58%  "

Set-TextValue "D17" "8.21"
$ws.Range("E17").Value = "  -8.60%  "

$ws.Range("D18").Value = "3.469.90"
$ws.Range("E18").Value = "  -4.30%  "

$ws.Range("D19").Value = "67.828.94"
$ws.Range("E19").Value = "  -3.93%  "

$ws.Range("E20").Value = "  -3.63%  "

$ws.Range("E21").Value = "  -3.88%  "

Set-TextValue "D22" "10.85"
$ws.Range("E22").Value = "  -4.86%  "

$ws.Range("E23").Value = "  -7.32%  "

Set-TextValue "D24" "15.39"
$ws.Range("E24").Value = "  -9.52%  "

Set-TextValue "D25" "94.81"
$ws.Range("E25").Value = "  -4.83%  "

$ws.Range("E26").Value = "  -5.10%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("E28").Value = "  -7.61%  "

Set-TextValue "D29" "8.93"
$ws.Range("E29").Value = "  -8.40%  "

Set-TextValue "D30" "31.94"
$ws.Range("E30").Value = "  -6.96%  "

Set-TextValue "D31" "8.28"
$ws.Range("E31").Value = "  -9.14%  "

Set-TextValue "D32" "3.02"
$ws.Range("E32").Value = "  -7.80%  "

$ws.Range("E33").Value = "  -7.16%  "

Set-TextValue "D34" "6.72"
$ws.Range("E34").Value = "  -9.52%  "

Set-TextValue "D35" "585.10"
$ws.Range("E35").Value = "  +2.59%  "

Set-TextValue "D36" "10.57"
$ws.Range("E36").Value = "  -4.49%  "

Set-TextValue "D37" "56.59"
$ws.Range("E37").Value = "  -3.15%  "

$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("E39").Value = "  -6.90%  "

Set-TextValue "D40" "3.34"
$ws.Range("E40").Value = "  -15.37%  "

$ws.Range("E41").Value = "  -3.83%  "

$ws.Range("E42").Value = "  -4.90%  "

$ws.Range("D43").Value = "3.357.37"
$ws.Range("E43").Value = "  -4.93%  "

$ws.Range("E44").Value = "  -8.20%  "

Set-TextValue "D45" "32.20"
$ws.Range("E45").Value = "  -5.80%  "

$ws.Range("D46").Value = "0.0₃0673"
$ws.Range("E46").Value = "  -7.17%  "

$ws.Range("E47").Value = "  -7.17%  "

$ws.Range("E48").Value = "  -9.13%  "

$ws.Range("E49").Value = "  -5.92%  "

Set-TextValue "D50" "131.90"
$ws.Range("E50").Value = "  -3.51%  "

Set-TextValue "D51" "5.55"
$ws.Range("E51").Value = "  +10.37%  "
